# Update scripts with new TPM values for the Col8a1-Itga1 LR-pairs sheet.
# Recomputed NATMI ligand/receptor expression, specificity and edge-weight
# columns (E,F,G,H,I,J,M,N,O,P,Q,R,S,T) for rows 2-10 to reflect the
# refreshed TPM inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 2.502470333333333
$ws.Cells.Item(2, 8).Value2 = 7.507410999999999
$ws.Cells.Item(2, 9).Value2 = 0.07821409705091072
$ws.Cells.Item(2, 10).Value2 = 0.07821409705091072
$ws.Cells.Item(2, 13).Value2 = 72.07569866666667
$ws.Cells.Item(2, 14).Value2 = 216.227096
$ws.Cells.Item(2, 15).Value2 = 0.4479522040449755
$ws.Cells.Item(2, 16).Value2 = 0.4479522040449755
$ws.Cells.Item(2, 17).Value2 = 180.3672976676062
$ws.Cells.Item(2, 18).Value2 = 1623.305679008456
$ws.Cells.Item(2, 19).Value2 = 0.03503617716134308
$ws.Cells.Item(2, 20).Value2 = 0.03503617716134308

$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 2.502470333333333
$ws.Cells.Item(3, 8).Value2 = 7.507410999999999
$ws.Cells.Item(3, 9).Value2 = 0.07821409705091072
$ws.Cells.Item(3, 10).Value2 = 0.07821409705091072
$ws.Cells.Item(3, 15).Value2 = 0.04737448730867841
$ws.Cells.Item(3, 16).Value2 = 0.0473744873086784
$ws.Cells.Item(3, 17).Value2 = 19.07526780110833
$ws.Cells.Item(3, 18).Value2 = 171.677410209975
$ws.Cells.Item(3, 19).Value2 = 0.003705352748098111
$ws.Cells.Item(3, 20).Value2 = 0.003705352748098111

$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 2.502470333333333
$ws.Cells.Item(4, 8).Value2 = 7.507410999999999
$ws.Cells.Item(4, 9).Value2 = 0.07821409705091072
$ws.Cells.Item(4, 10).Value2 = 0.07821409705091072
$ws.Cells.Item(4, 15).Value2 = 0.5046733086463462
$ws.Cells.Item(4, 16).Value2 = 0.5046733086463461
$ws.Cells.Item(4, 17).Value2 = 203.2059672070994
$ws.Cells.Item(4, 18).Value2 = 1828.853704863895
$ws.Cells.Item(4, 19).Value2 = 0.03947256714146954
$ws.Cells.Item(4, 20).Value2 = 0.03947256714146953

$ws.Cells.Item(5, 9).Value2 = 0.8193892102022395
$ws.Cells.Item(5, 10).Value2 = 0.8193892102022395
$ws.Cells.Item(5, 13).Value2 = 72.07569866666667
$ws.Cells.Item(5, 14).Value2 = 216.227096
$ws.Cells.Item(5, 15).Value2 = 0.4479522040449755
$ws.Cells.Item(5, 16).Value2 = 0.4479522040449755
$ws.Cells.Item(5, 17).Value2 = 1889.570079495679
$ws.Cells.Item(5, 18).Value2 = 17006.13071546111
$ws.Cells.Item(5, 19).Value2 = 0.367047202680765
$ws.Cells.Item(5, 20).Value2 = 0.367047202680765

$ws.Cells.Item(6, 9).Value2 = 0.8193892102022395
$ws.Cells.Item(6, 10).Value2 = 0.8193892102022395
$ws.Cells.Item(6, 15).Value2 = 0.04737448730867841
$ws.Cells.Item(6, 16).Value2 = 0.0473744873086784
$ws.Cells.Item(6, 19).Value2 = 0.03881814373959402
$ws.Cells.Item(6, 20).Value2 = 0.03881814373959402

$ws.Cells.Item(7, 9).Value2 = 0.8193892102022395
$ws.Cells.Item(7, 10).Value2 = 0.8193892102022395
$ws.Cells.Item(7, 15).Value2 = 0.5046733086463462
$ws.Cells.Item(7, 16).Value2 = 0.5046733086463461
$ws.Cells.Item(7, 19).Value2 = 0.4135238637818807
$ws.Cells.Item(7, 20).Value2 = 0.4135238637818806

$ws.Cells.Item(8, 7).Value2 = 3.276195666666666
$ws.Cells.Item(8, 8).Value2 = 9.828586999999999
$ws.Cells.Item(8, 9).Value2 = 0.1023966927468496
$ws.Cells.Item(8, 10).Value2 = 0.1023966927468496
$ws.Cells.Item(8, 13).Value2 = 72.07569866666667
$ws.Cells.Item(8, 14).Value2 = 216.227096
$ws.Cells.Item(8, 15).Value2 = 0.4479522040449755
$ws.Cells.Item(8, 16).Value2 = 0.4479522040449755
$ws.Cells.Item(8, 17).Value2 = 236.1340916437057
$ws.Cells.Item(8, 18).Value2 = 2125.206824793352
$ws.Cells.Item(8, 19).Value2 = 0.04586882420286746
$ws.Cells.Item(8, 20).Value2 = 0.04586882420286746

$ws.Cells.Item(9, 7).Value2 = 3.276195666666666
$ws.Cells.Item(9, 8).Value2 = 9.828586999999999
$ws.Cells.Item(9, 9).Value2 = 0.1023966927468496
$ws.Cells.Item(9, 10).Value2 = 0.1023966927468496
$ws.Cells.Item(9, 15).Value2 = 0.04737448730867841
$ws.Cells.Item(9, 16).Value2 = 0.0473744873086784
$ws.Cells.Item(9, 17).Value2 = 24.97304718384166
$ws.Cells.Item(9, 19).Value2 = 0.004850990820986271
$ws.Cells.Item(9, 20).Value2 = 0.004850990820986271

$ws.Cells.Item(10, 7).Value2 = 3.276195666666666
$ws.Cells.Item(10, 8).Value2 = 9.828586999999999
$ws.Cells.Item(10, 9).Value2 = 0.1023966927468496
$ws.Cells.Item(10, 10).Value2 = 0.1023966927468496
$ws.Cells.Item(10, 15).Value2 = 0.5046733086463462
$ws.Cells.Item(10, 16).Value2 = 0.5046733086463461
$ws.Cells.Item(10, 19).Value2 = 0.05167687772299593
$ws.Cells.Item(10, 20).Value2 = 0.05167687772299592
